$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new Price (column D) value. Values are stored as text
# (matching the source data format, which uses text cells for prices).
$priceUpdates = @{
    2 = "267.95"
    3 = "22.84"
    4 = "6.330"
    5 = "0.06198"
    8 = "1.396"
    9 = "0.8314"
    11 = "0.1611"
    12 = "0.08220"
    13 = "0.03399"
    14 = "0.03156"
    16 = "3.918"
    17 = "0.001724"
    18 = "0.04849"
    19 = "0.006312"
    20 = "0.005384"
    21 = "0.001090"
    22 = "0.0001501"
    24 = "2.366"
    26 = "0.1213"
    40 = "0.04655"
    41 = "0.006886"
    42 = "0.1153"
    43 = "0.003601"
    44 = "0.01214"
    45 = "0.00006264"
    48 = "0.1652"
    49 = "0.00002101"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item([int]$row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}
